# Commit: "Add new column 'Event' to Card24 by admin"
#
# The 'Event' column (M) header already exists on Card24; this edit clears
# the placeholder "nan" values that were populated in the data rows (M2:M12)
# for the newly added column, leaving those cells blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("M2:M12").ClearContents()
